$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.070500135421753
$ws.Range("B1").Value = 2.382978439331055
$ws.Range("C1").Value = 6.440464973449707
$ws.Range("D1").Value = 2.239310264587402
$ws.Range("E1").Value = 1.287963509559631
